$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: ANT / CHIP_ANTENNA (Mouser) ---
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "ANT"
$ws.Range("D11").Value = "CHIP_ANTENNA"
$ws.Range("E11").Value = "ANT"
$ws.Range("F11").Value = "2450AT18B100E"
$ws.Range("G11").Value = "Mouser"
$ws.Range("H11").Value = "609-2450AT18B100E"
$ws.Range("I11").Value = 5517
$ws.Range("I11").NumberFormat = "#,##0"
$ws.Range("J11").Value = 1.21
$url11 = "http://mx.mouser.com/search/ProductDetail.aspx?qs=yCnrNFeXz%252bh5MFsFIXGZGA==&utm_source=findchips&utm_medium=aggregator&utm_campaign=609-2450AT18B100E&utm_term=2450AT18B100"
$ws.Range("L11").Value = $url11
$ws.Hyperlinks.Add($ws.Range("L11"), $url11)
$ws.Range("M11").Value = "N/A"

# --- Row 12: ANT / CHIP_ANTENNA (Digikey) ---
$ws.Range("F12").Value = "2450AT18B100E"
$ws.Range("G12").Value = "Digikey"
$ws.Range("H12").Value = "712-1006-1-ND"
$ws.Range("I12").Value = 10828
$ws.Range("I12").NumberFormat = "#,##0"
$ws.Range("J12").Value = 0.94
$url12 = "http://www.digikey.com/product-detail/en/2450AT18B100E/712-1006-1-ND/1560835"
$ws.Range("L12").Value = $url12
$ws.Hyperlinks.Add($ws.Range("L12"), $url12)
$ws.Range("M12").Value = "N/A"

# --- Row 13: BAT165 / D1 (Mouser) ---
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "BAT165"
$ws.Range("D13").Value = "BAT165"
$ws.Range("E13").Value = "D1"
$ws.Range("F13").Value = "BAT 165 E6327"
$ws.Range("G13").Value = "Mouser"
$ws.Range("H13").Value = "726-BAT165E6327"
$ws.Range("I13").Value = 9000
$ws.Range("I13").NumberFormat = "#,##0"
$ws.Range("J13").Value = 0.702
$url13 = "http://mx.mouser.com/search/ProductDetail.aspx?qs=mzcOS1kGbgcQqWsJMFtrug==&utm_source=findchips&utm_medium=aggregator&utm_campaign=726-BAT165E6327&utm_term=BAT165"
$ws.Range("L13").Value = $url13
$ws.Hyperlinks.Add($ws.Range("L13"), $url13)
$ws.Range("M13").Value = "N/A"

# --- Row 14: BAT165 / D1 (Digikey) ---
$ws.Range("F14").Value = "BAT 165 E6327"
$ws.Range("G14").Value = "Digikey"
$ws.Range("H14").Value = "BAT 165 E6327CT-ND"
$ws.Range("I14").Value = 1619
$ws.Range("I14").NumberFormat = "#,##0"
$ws.Range("J14").Value = 0.55
$url14 = "http://www.digikey.com/product-detail/en/BAT%20165%20E6327/BAT%20165%20E6327CT-ND/3819504"
$ws.Range("L14").Value = $url14
$ws.Hyperlinks.Add($ws.Range("L14"), $url14)
$ws.Range("M14").Value = "N/A"

# --- Row 15: BALUM (Mouser) ---
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "BALUM"
$ws.Range("D15").Value = "BALUM"
$ws.Range("E15").Value = "BALUM"
$ws.Range("F15").Value = "2450BM14A0002T"
$ws.Range("G15").Value = "Mouser"
$ws.Range("H15").Value = "609-2450BM14A0002T"
$ws.Range("I15").Value = 3773
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("J15").Value = 1.91
$url15 = "http://mx.mouser.com/search/ProductDetail.aspx?qs=yCnrNFeXz%252bjc2NrpGmrycg==&utm_source=findchips&utm_medium=aggregator&utm_campaign=609-2450BM14A0002T&utm_term=2450BM14A0002"
$ws.Range("L15").Value = $url15
$ws.Hyperlinks.Add($ws.Range("L15"), $url15)
$ws.Range("M15").Value = "N/A"

# --- Row 16: SWITCH / SW (Mouser) ---
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = "SWITCH"
$ws.Range("D16").Value = "SWITCH"
$ws.Range("E16").Value = "SW"
$ws.Range("F16").Value = "SKQGADE010"
$ws.Range("G16").Value = "Mouser"
$ws.Range("H16").Value = "688-SKQGAD"
$ws.Range("I16").Value = 642
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("J16").Value = 0.923
$url16 = "http://mx.mouser.com/search/ProductDetail.aspx?qs=N5Jky1br14PCAY42dbciFw==&utm_source=findchips&utm_medium=aggregator&utm_campaign=688-SKQGAD&utm_term=SKQGADE010"
$ws.Range("L16").Value = $url16
$ws.Hyperlinks.Add($ws.Range("L16"), $url16)
$ws.Range("M16").Value = "N/A"

# Match the final selection left by the author's edit session
$ws.Range("J17").Select()
